$d = $word.ActiveDocument
$d.Content.Find.Execute("MySQL", $true, $false, $false, $false, $false, $true, 1, $false, "PostgreSQL", 2)
